$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("maximetro")

$ws.Range("C2").Value = 35
$ws.Range("D7").Value = 35
$ws.Range("D9").Value = 35
$ws.Range("D10").Value = 34
$ws.Range("C13").Value = 34
